# Femacal de La Calera - Papa: insert two new weekly-report rows.
#
# The source data block (rows 2-400) is a log of price reports ordered
# (mostly) oldest-to-newest per variety/quality group. This commit adds a
# new reporting date (2022-01-21, serial 44551) with two new "Rosara /
# 1a (cosecha)" observations at the very top of the "Rosara" block
# (rows 317-318), pushing every following row down by two positions
# (old row N -> new row N+2, for N = 317..400).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 317; this shifts the old
# rows 317..400 down to 319..402 and extends the sheet dimension
# automatically (A1:R400 -> A1:R402).
$ws.Rows.Item(317).Resize(2).Insert()

# --- New row 317 -----------------------------------------------------
$ws.Cells.Item(317, 1).Value = 3
$ws.Cells.Item(317, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(317, 3).Value = "Coquimbo"
$ws.Cells.Item(317, 4).Value = 44551
$ws.Cells.Item(317, 5).Value = 5
$ws.Cells.Item(317, 6).Value = 100114001
$ws.Cells.Item(317, 7).Value = "Papa"
$ws.Cells.Item(317, 8).Value = "Rosara"
$ws.Cells.Item(317, 9).Value = "1a (cosecha)"
$ws.Cells.Item(317, 10).Value = 450
$ws.Cells.Item(317, 11).Value = 8000
$ws.Cells.Item(317, 12).Value = 8500
$ws.Cells.Item(317, 13).Value = 8222
$ws.Cells.Item(317, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(317, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(317, 16).Value = 329
$ws.Cells.Item(317, 17).Value = 25
$ws.Cells.Item(317, 18).Value = "Hortaliza"

# --- New row 318 -----------------------------------------------------
$ws.Cells.Item(318, 1).Value = 3
$ws.Cells.Item(318, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(318, 3).Value = "Coquimbo"
$ws.Cells.Item(318, 4).Value = 44551
$ws.Cells.Item(318, 5).Value = 5
$ws.Cells.Item(318, 6).Value = 100114001
$ws.Cells.Item(318, 7).Value = "Papa"
$ws.Cells.Item(318, 8).Value = "Rosara"
$ws.Cells.Item(318, 9).Value = "1a (cosecha)"
$ws.Cells.Item(318, 10).Value = 570
$ws.Cells.Item(318, 11).Value = 8500
$ws.Cells.Item(318, 12).Value = 9000
$ws.Cells.Item(318, 13).Value = 8754
$ws.Cells.Item(318, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(318, 15).Value = "Provincia de Talca"
$ws.Cells.Item(318, 16).Value = 350
$ws.Cells.Item(318, 17).Value = 25
$ws.Cells.Item(318, 18).Value = "Hortaliza"

# D column (Fecha) keeps the date-number-format style already present
# on the rest of the column; re-apply it explicitly in case Insert()
# didn't carry it through for the two brand-new cells.
$ws.Range("D317:D318").NumberFormat = "YYYY-MM-DD HH:MM:SS"
